$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 167162
$ws.Range("C4").Value = 158062
$ws.Range("C7").Value = 5.44
$ws.Range("C8").Value = 65.29
